# Append a new row (row 86) to the single worksheet "Sheet1", extending
# the existing games table: Nome do Jogo | Status | Plataforma | Objetivo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the numeric-looking value "5", but (like the rows just
# above it, A83:A85 = "1"/"2"/"3") it must be stored as TEXT, not as a
# number. A plain Value assignment of "5" would be auto-coerced to a
# number by Excel, so instead we write a formula that evaluates to the
# text string "5" and then convert it to a literal via copy / paste
# values - this yields a genuine text cell without leaving behind any
# stray "quote prefix" cell formatting.
$ws.Range("A86").Formula = "=""5"""
$ws.Range("A86").Copy()
$ws.Range("A86").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B86").Value = "Completo"
$ws.Range("C86").Value = "PS3"
$ws.Range("D86").Value = "Platinado"
